$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H3").Value = 125657
$ws.Range("J3").Value = 125657
$ws.Range("L3").Value = 125657
$ws.Range("N3").Value = -125885
$ws.Range("H11").Value = 92.7
$ws.Range("I11").Value = 92.7
$ws.Range("K11").Value = 92.7
$ws.Range("M11").Value = 47.3
$ws.Range("H70").Value = 3054.8
$ws.Range("I70").Value = 3599.6
$ws.Range("J70").Value = 2510
$ws.Range("K70").Value = 10798.8
$ws.Range("L70").Value = 7530
$ws.Range("M70").Value = -10528.8
$ws.Range("N70").Value = -8070
$ws.Range("H73").Value = 3054.8
$ws.Range("I73").Value = 3599.6
$ws.Range("J73").Value = 2510
$ws.Range("K73").Value = 10798.8
$ws.Range("L73").Value = 7530
$ws.Range("M73").Value = -9862.799999999999
$ws.Range("N73").Value = -9402
$ws.Range("H87").Value = 54997.5
$ws.Range("J87").Value = 54997.5
$ws.Range("L87").Value = 54997.5
$ws.Range("N87").Value = -57493.5
$ws.Range("H88").Value = 99999
$ws.Range("I88").Value = 0
$ws.Range("J88").Value = 99999
$ws.Range("K88").Value = 0
$ws.Range("L88").Value = 99999
$ws.Range("M88").ClearContents()
$ws.Range("N88").Value = -100811
$ws.Range("H90").Value = 54997.5
$ws.Range("J90").Value = 54997.5
$ws.Range("L90").Value = 164992.5
$ws.Range("N90").Value = -177472.5
$ws.Range("H91").Value = 99999
$ws.Range("I91").Value = 0
$ws.Range("J91").Value = 99999
$ws.Range("K91").Value = 0
$ws.Range("L91").Value = 99999
$ws.Range("M91").ClearContents()
$ws.Range("N91").Value = -102807
$ws.Range("H102").Value = 125657
$ws.Range("J102").Value = 125657
$ws.Range("L102").Value = 125657
$ws.Range("N102").Value = -132147
$ws.Range("H113").Value = 6477.357
$ws.Range("I113").Value = 5250
$ws.Range("J113").Value = 6968.3
$ws.Range("K113").Value = 5250
$ws.Range("L113").Value = 6968.3
$ws.Range("M113").Value = -1996
$ws.Range("N113").Value = -13476.3
$ws.Range("H116").Value = 4331.6665
$ws.Range("I116").Value = 3997.5
$ws.Range("K116").Value = 3997.5
$ws.Range("M116").Value = -555.5
$ws.Range("H132").Value = 1097.3636
$ws.Range("I132").Value = 857.6
$ws.Range("K132").Value = 2572.8
$ws.Range("M132").Value = -42.80000000000018

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 2285.3845
$ws.Range("I61").Value = 2246.3635
$ws.Range("J61").Value = 2500
$ws.Range("K61").Value = 2246.3635
$ws.Range("L61").Value = 2500
$ws.Range("M61").Value = -2034.3635
$ws.Range("N61").Value = -2924
$ws.Range("H136").Value = 2285.3845
$ws.Range("I136").Value = 2246.3635
$ws.Range("J136").Value = 2500
$ws.Range("K136").Value = 6739.0905
$ws.Range("L136").Value = 7500
$ws.Range("M136").Value = -4189.0905
$ws.Range("N136").Value = -12600

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H86").Value = 11166.667
$ws.Range("I86").Value = 10000
$ws.Range("J86").Value = 11400
$ws.Range("K86").Value = 10000
$ws.Range("L86").Value = 11400
$ws.Range("M86").Value = -8877
$ws.Range("N86").Value = -13646
$ws.Range("H89").Value = 11166.667
$ws.Range("I89").Value = 10000
$ws.Range("J89").Value = 11400
$ws.Range("K89").Value = 50000
$ws.Range("L89").Value = 57000
$ws.Range("M89").Value = -44384
$ws.Range("N89").Value = -68232
$ws.Range("H132").Value = 2746.3635
$ws.Range("I132").Value = 2280.75
$ws.Range("K132").Value = 6842.25
$ws.Range("M132").Value = -4312.25
$ws.Range("H141").Value = 21000
$ws.Range("J141").Value = 21000
$ws.Range("L141").Value = 21000
$ws.Range("N141").Value = -31360

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H97").Value = 1989.75
$ws.Range("I97").Value = 1345.4286
$ws.Range("K97").Value = 1345.4286
$ws.Range("M97").Value = -849.4286
$ws.Range("H101").Value = 63119
$ws.Range("J101").Value = 63119
$ws.Range("L101").Value = 63119
$ws.Range("N101").Value = -69609
$ws.Range("H102").Value = 13972.6
$ws.Range("I102").Value = 3290.111
$ws.Range("K102").Value = 3290.111
$ws.Range("M102").Value = -1668.111
$ws.Range("H126").Value = 15000
$ws.Range("I126").Value = 15000
$ws.Range("J126").Value = 0
$ws.Range("K126").Value = 45000
$ws.Range("L126").Value = 0
$ws.Range("M126").Value = -42530
$ws.Range("N126").ClearContents()

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 20674.23
$ws.Range("I7").Value = 20674.23
$ws.Range("K7").Value = 20674.23
$ws.Range("M7").Value = -20562.23
$ws.Range("H61").Value = 5000
$ws.Range("I61").Value = 5000
$ws.Range("K61").Value = 5000
$ws.Range("M61").Value = -4798
$ws.Range("H63").Value = 43028.332
$ws.Range("J63").Value = 43028.332
$ws.Range("L63").Value = 43028.332
$ws.Range("N63").Value = -44526.332
$ws.Range("H66").Value = 43028.332
$ws.Range("J66").Value = 43028.332
$ws.Range("L66").Value = 129084.996
$ws.Range("N66").Value = -136572.996
$ws.Range("H113").Value = 5000
$ws.Range("I113").Value = 5000
$ws.Range("K113").Value = 5000
$ws.Range("M113").Value = -2830
$ws.Range("H122").Value = 3420.6667
$ws.Range("I122").Value = 3403.8
$ws.Range("K122").Value = 10211.4
$ws.Range("M122").Value = -7761.400000000001
$ws.Range("H126").Value = 20674.23
$ws.Range("I126").Value = 20674.23
$ws.Range("K126").Value = 62022.69
$ws.Range("M126").Value = -59552.69

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H81").Value = 994
$ws.Range("I81").Value = 994
$ws.Range("K81").Value = 1988
$ws.Range("M81").Value = -927
$ws.Range("H84").Value = 994
$ws.Range("I84").Value = 994
$ws.Range("K84").Value = 9940
$ws.Range("M84").Value = -4636
$ws.Range("H107").Value = 1341
$ws.Range("I107").Value = 893.875
$ws.Range("J107").Value = 2533.3333
$ws.Range("K107").Value = 2681.625
$ws.Range("L107").Value = 7599.999899999999
$ws.Range("M107").Value = -761.625
$ws.Range("N107").Value = -11439.9999
